$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 544 (shifts existing rows 544.. down by one)
$ws.Rows.Item(544).Insert()

# Populate the newly inserted row 544 with the new price record
$ws.Cells.Item(544, 1).Value  = 9
$ws.Cells.Item(544, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(544, 3).Value  = "Metropolitana"
$ws.Cells.Item(544, 4).Value  = 45218
$ws.Cells.Item(544, 5).Value  = 13
$ws.Cells.Item(544, 6).Value  = 100112044
$ws.Cells.Item(544, 7).Value  = "Perejil"
$ws.Cells.Item(544, 8).Value  = "Sin especificar"
$ws.Cells.Item(544, 9).Value  = "Primera"
$ws.Cells.Item(544, 10).Value = 70
$ws.Cells.Item(544, 11).Value = 13000
$ws.Cells.Item(544, 12).Value = 14000
$ws.Cells.Item(544, 13).Value = 13500
$ws.Cells.Item(544, 14).Value = "$/docena de atados"
$ws.Cells.Item(544, 15).Value = "Región Metropolitana"
$ws.Cells.Item(544, 16).Value = 4500
$ws.Cells.Item(544, 17).Value = 3
$ws.Cells.Item(544, 18).Value = "Hortaliza"

# Make sure the date column keeps its date style (numFmtId 165), matching the
# formatting used by the rest of column D.
$ws.Cells.Item(544, 4).NumberFormat = $ws.Cells.Item(545, 4).NumberFormat
